$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues = -4163

# ---------------------------------------------------------------------------
# Row 32: "3.6" label loses its gold highlight fill, and its description is
# replaced by the new (more detailed) text. A new "x" mark is added in K32.
# ---------------------------------------------------------------------------

# Remove the fill/highlight from B32 ("3.6") by copying the plain bold/centered
# format already used for the other "3.x" labels (B30), leaving its value as is.
$ws.Range("B30").Copy()
$ws.Range("B32").PasteSpecial($xlPasteFormats)

$ws.Range("C32").Value = "Experiencia docente cursos previos promedio por dpto (min, max, avg, stddev)"

$ws.Range("K32").Value = "x"

# ---------------------------------------------------------------------------
# Row 33 (new): "3.7" label + its description.
# ---------------------------------------------------------------------------

# "3.7" needs to be stored as text (like the other "3.x" labels), not a number.
# Build it as a text-formula result in a scratch cell far outside the used
# range, bring only the resulting value over (so it lands as a plain shared
# string), then restore the normal label formatting and clean up the scratch
# cell, leaving no left-over formatting artifacts in the workbook.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '="3.7"'
$scratch.Copy()
$ws.Range("B33").PasteSpecial($xlPasteValues)
$scratch.Clear()

$ws.Range("B30").Copy()
$ws.Range("B33").PasteSpecial($xlPasteFormats)

$ws.Range("C33").Value = "Experiencia docente cursos que va a ver en el periodo actual promedio por dpto (min, max, avg, stddev)"

$ws.Rows("33").RowHeight = 68

# ---------------------------------------------------------------------------
# Row 41 (new): FECHA 1 / 9 DE DICIEMBRE
# ---------------------------------------------------------------------------
$ws.Range("B41").Value = "FECHA 1"
$ws.Range("C41").Value = "9 DE DICIEMBRE"
$ws.Rows("41").RowHeight = 17

# ---------------------------------------------------------------------------
# Row 42 (new): FECHA 2 / 22 DE ENERO (date-styled text)
# ---------------------------------------------------------------------------
$ws.Range("B42").Value = "FECHA 2"
$ws.Range("C42").Value = "22 DE ENERO"
$ws.Range("C42").NumberFormat = "d-mmm"
$ws.Rows("42").RowHeight = 17

# ---------------------------------------------------------------------------
# Row 45-48 (new): VIDEO notes
# ---------------------------------------------------------------------------
$ws.Range("B45").Value = "VIDEO"
$ws.Range("C45").Value = "ya usar un powerpoint de cómo voy "
$ws.Rows("45").RowHeight = 34

$ws.Range("C46").Value = "dedicarse solo al tema de procesamiento, modelos"
$ws.Rows("46").RowHeight = 34

$ws.Range("C47").Value = "no el sustento del problema etc…."
$ws.Rows("47").RowHeight = 34

$ws.Range("C48").Value = "FECHA:  3 DE NOVIEMBRE"
$ws.Rows("48").RowHeight = 17

# ---------------------------------------------------------------------------
# View: zoom out a bit and move the selection to C32.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 130
$ws.Range("C32").Select()
